$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (row 1, A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 13:25"

# --- Update per-country statistics ---
# Row 16: Austria
$ws.Range("B16").Value = 8450
$ws.Range("C16").Value = 179
$ws.Range("E16").Value = 7885

# Row 20: Noruega
$ws.Range("B20").Value = 4189
$ws.Range("C20").Value = 174
$ws.Range("E20").Value = 4159

# Row 25: Chequia
$ws.Range("B25").Value = 2689
$ws.Range("C25").Value = 58
$ws.Range("E25").Value = 2665

# Row 28: Dinamarca
$ws.Range("B28").Value = 2395
$ws.Range("C28").Value = 194
$ws.Range("E28").Value = 2322
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 72

# Row 32: Rumania
$ws.Range("E32").Value = 1551
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 40

# Row 78: Taiwan
$ws.Range("D78").Value = 39
$ws.Range("E78").Value = 257

# Rows 81-84: a new entry for "Republica de Macedonia" is inserted ahead of
# Kuwait, Kazajistan and Jordania, so those three shift down one row.
# Row 81: Republica de Macedonia (new data)
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 257
$ws.Range("C81").Value = 16
$ws.Range("D81").Value = 3
$ws.Range("E81").Value = 248
$ws.Range("F81").Value = 1
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 6

# Row 82: Kuwait (shifted down, values unchanged from former row 81)
$ws.Range("A82").Value = "Kuwait"
$ws.Range("B82").Value = 255
$ws.Range("C82").Value = 20
$ws.Range("D82").Value = 67
$ws.Range("E82").Value = 188
$ws.Range("F82").Value = 12
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 0

# Row 83: Kazajistan (shifted down, values unchanged from former row 82)
$ws.Range("A83").Value = "Kazajistan"
$ws.Range("B83").Value = 251
$ws.Range("C83").Value = 23
$ws.Range("D83").Value = 18
$ws.Range("E83").Value = 232
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 1

# Row 84: Jordania (shifted down, values unchanged from former row 83)
$ws.Range("A84").Value = "Jordania"
$ws.Range("B84").Value = 246
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 18
$ws.Range("E84").Value = 227
$ws.Range("F84").Value = 3
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 1

# Row 95: Malta
$ws.Range("B95").Value = 151
$ws.Range("C95").Value = 2
$ws.Range("E95").Value = 149
